$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (col C and D narrower) ---
# Target OOXML widths: C=10.1530612244898, D=10.0051020408163
# Engine snaps ColumnWidth to an internal pixel grid (6px/char), so the
# inputs below are the closest achievable values.
$ws.Columns.Item(3).ColumnWidth = 9.333333333333334
$ws.Columns.Item(4).ColumnWidth = 9.166666666666666

# --- View: scroll position + active selection ---
$ws.Range("A7").Select() | Out-Null
$ws.Range("F28").Select() | Out-Null

# --- New AVERAGE formulas in column F (summary rows) ---
$ws.Range("F13").Formula = "=AVERAGE(E10:E16)"
$ws.Range("F22").Formula = "=AVERAGE(E18:E24)"
$ws.Range("F27").Formula = "=AVERAGE(E26:E32)"

# Row 47
$ws.Range("Q47").Value = "predict"
$ws.Range("R47").Value = "mispred"

# Row 48
$ws.Range("N48").Value = "ackermann"

# Row 49
$ws.Range("O49").Value = "Register Mapped"
$ws.Range("P49").Value = 31492889979
$ws.Range("Q49").Formula = "=P49-R49"
$ws.Range("R49").Value = 2863442739
$ws.Range("S49").Value = 1

# Row 50
$ws.Range("O50").Value = "Normal"
$ws.Range("P50").Value = 31492889979
$ws.Range("Q50").Formula = "=P50-R50"
$ws.Range("R50").Value = 17177904179
$ws.Range("S50").Value = 3

# Row 51
$ws.Range("O51").Value = "Table Only"
$ws.Range("P51").Value = 31492889979
$ws.Range("Q51").Formula = "=P51-R51"
$ws.Range("R51").Value = 2863442739
$ws.Range("S51").Value = 2

# Row 52
$ws.Range("N52").Value = "fasta"

# Row 53
$ws.Range("O53").Value = "Register Mapped"
$ws.Range("P53").Value = 10202767734
$ws.Range("Q53").Formula = "=P53-R53"
$ws.Range("R53").Value = 2923312485
$ws.Range("S53").Value = 1

# Row 54
$ws.Range("O54").Value = "Normal"
$ws.Range("P54").Value = 10202767734
$ws.Range("Q54").Formula = "=P54-R54"
$ws.Range("R54").Value = 5157715191
$ws.Range("S54").Value = 3

# Row 55
$ws.Range("O55").Value = "Table Only"
$ws.Range("P55").Value = 10202767734
$ws.Range("Q55").Formula = "=P55-R55"
$ws.Range("R55").Value = 3520063385
$ws.Range("S55").Value = 2

# Row 56
$ws.Range("N56").Value = "reversecomplement"

# Row 57
$ws.Range("O57").Value = "Register Mapped"
$ws.Range("P57").Value = 3079415632
$ws.Range("Q57").Formula = "=P57-R57"
$ws.Range("R57").Value = 302084202
$ws.Range("S57").Value = 1

# Row 58
$ws.Range("O58").Value = "Normal"
$ws.Range("P58").Value = 3079415632
$ws.Range("Q58").Formula = "=P58-R58"
$ws.Range("R58").Value = 552207786
$ws.Range("S58").Value = 3

# Row 59
$ws.Range("O59").Value = "Table Only"
$ws.Range("P59").Value = 3079415632
$ws.Range("Q59").Formula = "=P59-R59"
$ws.Range("R59").Value = 54166945
$ws.Range("S59").Value = 2

# Row 60
$ws.Range("N60").Value = "mersenne"

# Row 61
$ws.Range("O61").Value = "Register Mapped"
$ws.Range("P61").Value = 6401505120
$ws.Range("Q61").Formula = "=P61-R61"
$ws.Range("R61").Value = 1851826825
$ws.Range("S61").Value = 1

# Row 62
$ws.Range("O62").Value = "Normal"
$ws.Range("P62").Value = 6401505120
$ws.Range("Q62").Formula = "=P62-R62"
$ws.Range("R62").Value = 3301361969
$ws.Range("S62").Value = 3

# Row 63
$ws.Range("O63").Value = "Table Only"
$ws.Range("P63").Value = 6401505120
$ws.Range("Q63").Formula = "=P63-R63"
$ws.Range("R63").Value = 1502144873
$ws.Range("S63").Value = 2

# Row 64
$ws.Range("N64").Value = "fannkuch"

# Row 65
$ws.Range("O65").Value = "Register Mapped"
$ws.Range("P65").Value = 13506897459
$ws.Range("Q65").Formula = "=P65-R65"
$ws.Range("R65").Value = 3069689230
$ws.Range("S65").Value = 1

# Row 66
$ws.Range("O66").Value = "Normal"
$ws.Range("P66").Value = 13506897459
$ws.Range("Q66").Formula = "=P66-R66"
$ws.Range("R66").Value = 8989685254
$ws.Range("S66").Value = 3

# Row 67
$ws.Range("O67").Value = "Table Only"
$ws.Range("P67").Value = 13506897459
$ws.Range("Q67").Formula = "=P67-R67"
$ws.Range("R67").Value = 4262582328
$ws.Range("S67").Value = 2

# Row 68
$ws.Range("N68").Value = "primesieve"

# Row 69
$ws.Range("O69").Value = "Register Mapped"
$ws.Range("P69").Value = 24690256810
$ws.Range("Q69").Formula = "=P69-R69"
$ws.Range("R69").Value = 3314133
$ws.Range("S69").Value = 1

# Row 70
$ws.Range("O70").Value = "Normal"
$ws.Range("P70").Value = 24690256810
$ws.Range("Q70").Formula = "=P70-R70"
$ws.Range("R70").Value = 12343774987
$ws.Range("S70").Value = 3

# Row 71
$ws.Range("O71").Value = "Table Only"
$ws.Range("P71").Value = 24690256810
$ws.Range("Q71").Formula = "=P71-R71"
$ws.Range("R71").Value = 3157143
$ws.Range("S71").Value = 2

# Row 72
$ws.Range("N72").Value = "mandelbrot"

# Row 73
$ws.Range("O73").Value = "Register Mapped"
$ws.Range("P73").Value = 182287976534
$ws.Range("Q73").Formula = "=P73-R73"
$ws.Range("R73").Value = 54104775594
$ws.Range("S73").Value = 1

# Row 74
$ws.Range("O74").Value = "Normal"
$ws.Range("P74").Value = 182287976534
$ws.Range("Q74").Formula = "=P74-R74"
$ws.Range("R74").Value = 107030190090
$ws.Range("S74").Value = 3

# Row 75
$ws.Range("O75").Value = "Table Only"
$ws.Range("P75").Value = 182287976534
$ws.Range("Q75").Formula = "=P75-R75"
$ws.Range("R75").Value = 54104759595
$ws.Range("S75").Value = 2

Write-Output "edit complete"